# Apply the citation-check / document-browser update described by the diff:
#   1. Strip the Heading2 paragraph style from the four section headings
#      (Introduction, Geographic Area and Needs Assessment,
#       Health Equity and Proposed Changes, Conclusion) so they revert
#      to the default (Normal) paragraph style.
#   2. Replace the inline author/year citation markers with the new
#      "Ref-XXXXXXX" placeholders (or "Smith, 2021" for the final body
#      paragraph), scoping every Find/Replace to its own paragraph so
#      that citations repeated verbatim elsewhere in the document are
#      left untouched.

$d = $word.ActiveDocument

function Replace-InParagraph($paraIndex, $oldText, $newText) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $r.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
}

# -- 1. De-style the four Heading2 paragraphs -----------------------------
$headingIndexes = @(2, 6, 12, 20)
foreach ($idx in $headingIndexes) {
    $d.Paragraphs($idx).Style = "Normal"
}

# -- 2. Citation swaps, scoped paragraph-by-paragraph ----------------------

# Paragraph 8 - "Geographic Area and Needs Assessment" body
Replace-InParagraph 8 "(Nuako et al.)" "(Ref-f470498)"
Replace-InParagraph 8 "(Cole)" "(Ref-f470498)"

# Paragraph 10 - "Furthermore, adherence to governmental regulations..."
Replace-InParagraph 10 "(Cole)" "(Ref-s908197)"
Replace-InParagraph 10 "(Johnston et al.)" "(Ref-s908197)"

# Paragraph 14 - "To advance health equity..."
Replace-InParagraph 14 "(Lyles et al.)" "(Ref-s595245)"
Replace-InParagraph 14 "(Patel et al.)" "(Ref-s595245)"

# Paragraph 16 - "Additionally, a horizon scan..."
Replace-InParagraph 16 "(Cole)" "(Ref-u951518)"
Replace-InParagraph 16 "(Johnston et al.)" "(Ref-u951518)"

# Paragraph 18 - "Moreover, the proposal emphasizes..."
Replace-InParagraph 18 "(Nuako et al.)" "(Smith, 2021)"
Replace-InParagraph 18 "(Lyles et al.)" "(Smith, 2021)"
